$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "durum" sheet: StartDate column (B2:B7) moves from 2025-06-24 -> 2025-06-26
#    Use a throwaway formula + paste-values round-trip so Excel's automatic
#    date detection doesn't turn the literal string into a date serial.
# ---------------------------------------------------------------------------
$durum = $wb.Worksheets.Item("durum")
$durum.Range("B2:B7").Formula = '="2025-06-26"'
$durum.Range("B2:B7").Copy()
$durum.Range("B2:B7").PasteSpecial(-4163)

# ---------------------------------------------------------------------------
# 2) District detail sheets: append one new row (new price date 45834) at the
#    bottom of each sheet's table (rows 2-197 -> 2-198).
# ---------------------------------------------------------------------------
$districtNewRows = @{
    "934015" = 49.06
    "065001" = 51.31
    "035001" = 50.26
    "055001" = 50.54
    "021001" = 51.26
    "038001" = 50.94
}

foreach ($sheetName in $districtNewRows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 198
    $ws.Cells.Item($newRow, 1).Value = 45834
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 2).Value = "Motorin UltraForce"
    $ws.Cells.Item($newRow, 3).Value = $districtNewRows[$sheetName]
}

# ---------------------------------------------------------------------------
# 3) "eskalasyon" sheet: each of the 9 "Name" groups gains a new trailing row
#    for price date 45834 (rate stays 0.05, amount/degisim/eskalasyon below).
#    Insert bottom-up so earlier row numbers stay valid while we work.
# ---------------------------------------------------------------------------
$esk = $wb.Worksheets.Item("eskalasyon")

# (insertAfterRow, amount, degisim/eskalasyon)
$groups = @(
    @{ After = 45; Amount = 49.06; Change = -0.1032717967464815 },  # TL/Desi Avrupa&Anadolu Dağıtım
    @{ After = 40; Amount = 49.06; Change = -0.1032717967464815 },  # TL/Desi Avrupa İade Toplama
    @{ After = 36; Amount = 49.06; Change = -0.1032717967464815 },  # TL/Desi Avrupa Toplama
    @{ After = 31; Amount = 49.06; Change = -0.1032717967464815 },  # Spot Araç Teknosa
    @{ After = 26; Amount = 49.06; Change = -0.1032717967464815 },  # Spot Araç Avrupa&Anadolu
    @{ After = 21; Amount = 49.06; Change = -0.1032717967464815 },  # Spot Araç Anadolu Toplama
    @{ After = 16; Amount = 50.54; Change = -0.1005516995906744 },  # Servis Samsun
    @{ After = 11; Amount = 50.94; Change = -0.09984096130058329 }, # Servis Kayseri
    @{ After = 6;  Amount = 51.26; Change = -0.09927956422421369 } # Servis Diyarbakır
)

foreach ($g in $groups) {
    $targetRow = $g.After + 1
    if ($targetRow -le 45) {
        $esk.Rows.Item($targetRow).Insert()
    }
    $nameVal = $esk.Cells.Item($g.After, 6).Value2
    $esk.Cells.Item($targetRow, 1).Value = 45834
    $esk.Cells.Item($targetRow, 2).Value = "Motorin UltraForce"
    $esk.Cells.Item($targetRow, 3).Value = $g.Amount
    $esk.Cells.Item($targetRow, 4).Value = $g.Change
    $esk.Cells.Item($targetRow, 5).Value = $g.Change
    $esk.Cells.Item($targetRow, 6).Value = $nameVal
    $esk.Cells.Item($targetRow, 7).Value = 0.05
}
